$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ccl2"
$ws.Range("C2").Value = "Ccr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 56.47185866666666
$ws.Range("H2").Value = 169.415576
$ws.Range("I2").Value = 0.8070274173741353
$ws.Range("J2").Value = 0.8070274173741354
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 9.235433666666667
$ws.Range("N2").Value = 27.706301
$ws.Range("O2").Value = 0.9849159412561933
$ws.Range("P2").Value = 0.9849159412561931
$ws.Range("Q2").Value = 521.5421047493751
$ws.Range("R2").Value = 4693.878942744376
$ws.Range("S2").Value = 0.7948541684026013
$ws.Range("T2").Value = 0.7948541684026011

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ccl2"
$ws.Range("C3").Value = "Ccr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 56.47185866666666
$ws.Range("H3").Value = 169.415576
$ws.Range("I3").Value = 0.8070274173741353
$ws.Range("J3").Value = 0.8070274173741354
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.08345333333333332
$ws.Range("N3").Value = 0.25036
$ws.Range("O3").Value = 0.008899908907107467
$ws.Range("P3").Value = 0.008899908907107465
$ws.Range("Q3").Value = 4.712764845262221
$ws.Range("R3").Value = 42.41488360735999
$ws.Range("S3").Value = 0.007182470500168002
$ws.Range("T3").Value = 0.007182470500168001

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ccl2"
$ws.Range("C4").Value = "Ccr2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 56.47185866666666
$ws.Range("H4").Value = 169.415576
$ws.Range("I4").Value = 0.8070274173741353
$ws.Range("J4").Value = 0.8070274173741354
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.057988
$ws.Range("N4").Value = 0.173964
$ws.Range("O4").Value = 0.006184149836699327
$ws.Range("P4").Value = 0.006184149836699326
$ws.Range("Q4").Value = 3.274690140362666
$ws.Range("R4").Value = 29.472211263264
$ws.Range("S4").Value = 0.004990778471366138
$ws.Range("T4").Value = 0.004990778471366138

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ccl2"
$ws.Range("C5").Value = "Ccr2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 10.23495333333333
$ws.Range("H5").Value = 30.70486
$ws.Range("I5").Value = 0.1462655586439962
$ws.Range("J5").Value = 0.1462655586439962
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 9.235433666666667
$ws.Range("N5").Value = 27.706301
$ws.Range("O5").Value = 0.9849159412561933
$ws.Range("P5").Value = 0.9849159412561931
$ws.Range("Q5").Value = 94.5242325914289
$ws.Range("R5").Value = 850.7180933228601
$ws.Range("S5").Value = 0.1440592803652144
$ws.Range("T5").Value = 0.1440592803652144

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Ccl2"
$ws.Range("C6").Value = "Ccr2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 10.23495333333333
$ws.Range("H6").Value = 30.70486
$ws.Range("I6").Value = 0.1462655586439962
$ws.Range("J6").Value = 0.1462655586439962
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.08345333333333332
$ws.Range("N6").Value = 0.25036
$ws.Range("O6").Value = 0.008899908907107467
$ws.Range("P6").Value = 0.008899908907107465
$ws.Range("Q6").Value = 0.8541409721777777
$ws.Range("R6").Value = 7.687268749599999
$ws.Range("S6").Value = 0.001301750148178751
$ws.Range("T6").Value = 0.001301750148178751

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ccl2"
$ws.Range("C7").Value = "Ccr2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 10.23495333333333
$ws.Range("H7").Value = 30.70486
$ws.Range("I7").Value = 0.1462655586439962
$ws.Range("J7").Value = 0.1462655586439962
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.057988
$ws.Range("N7").Value = 0.173964
$ws.Range("O7").Value = 0.006184149836699327
$ws.Range("P7").Value = 0.006184149836699326
$ws.Range("Q7").Value = 0.5935044738933334
$ws.Range("R7").Value = 5.34154026504
$ws.Range("S7").Value = 0.0009045281306030048
$ws.Range("T7").Value = 0.0009045281306030048

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Ccl2"
$ws.Range("C8").Value = "Ccr2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3.268330666666667
$ws.Range("H8").Value = 9.804992
$ws.Range("I8").Value = 0.04670702398186845
$ws.Range("J8").Value = 0.04670702398186846
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 9.235433666666667
$ws.Range("N8").Value = 27.706301
$ws.Range("O8").Value = 0.9849159412561933
$ws.Range("P8").Value = 0.9849159412561931
$ws.Range("Q8").Value = 30.18445107273245
$ws.Range("R8").Value = 271.660059654592
$ws.Range("S8").Value = 0.04600249248837756
$ws.Range("T8").Value = 0.04600249248837756

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Ccl2"
$ws.Range("C9").Value = "Ccr2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3.268330666666667
$ws.Range("H9").Value = 9.804992
$ws.Range("I9").Value = 0.04670702398186845
$ws.Range("J9").Value = 0.04670702398186846
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.08345333333333332
$ws.Range("N9").Value = 0.25036
$ws.Range("O9").Value = 0.008899908907107467
$ws.Range("P9").Value = 0.008899908907107465
$ws.Range("Q9").Value = 0.2727530885688889
$ws.Range("R9").Value = 2.45477779712
$ws.Range("S9").Value = 0.0004156882587607132
$ws.Range("T9").Value = 0.0004156882587607131

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Ccl2"
$ws.Range("C10").Value = "Ccr2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 3.268330666666667
$ws.Range("H10").Value = 9.804992
$ws.Range("I10").Value = 0.04670702398186845
$ws.Range("J10").Value = 0.04670702398186846
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.057988
$ws.Range("N10").Value = 0.173964
$ws.Range("O10").Value = 0.006184149836699327
$ws.Range("P10").Value = 0.006184149836699326
$ws.Range("Q10").Value = 0.1895239586986667
$ws.Range("R10").Value = 1.705715628288
$ws.Range("S10").Value = 0.0002888432347301834
$ws.Range("T10").Value = 0.0002888432347301834

